$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New column D: thresh_expected - boolean results for low/high threshold fuzzy-match tests.
# Rows where score falls strictly between the low and high thresholds are left blank (None).
$ws.Range("D1").Value = "thresh_expected"

$ws.Range("D2").Value = $true
$ws.Range("D3").Value = $true
$ws.Range("D4").Value = $true
$ws.Range("D5").Value = $true
$ws.Range("D6").Value = $true
$ws.Range("D7").Value = $true
$ws.Range("D8").Value = $true
$ws.Range("D9").Value = $true
$ws.Range("D10").Value = $false
$ws.Range("D11").Value = $false
# D12:D15 intentionally left blank (None - ambiguous low < score < high)
$ws.Range("D16").Value = $false
$ws.Range("D17").Value = $false

$ws.Range("D2:D17").Select() | Out-Null
